# Swap the contents of columns B..G between paired rows.
# (Each pair of adjacent stock-report lines had its data rows
#  transposed; column A, the running serial number, stays put.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(33, 34),
    @(151, 152),
    @(198, 199),
    @(228, 229),
    @(237, 238),
    @(326, 327),
    @(371, 372),
    @(373, 374),
    @(387, 388),
    @(391, 392),
    @(401, 402),
    @(484, 485),
    @(560, 561),
    @(573, 574),
    @(673, 674),
    @(834, 835)
)

$cols = @(2, 3, 4, 5, 6, 7)  # B, C, D, E, F, G

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
